$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.729.94"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.851.94"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  -2.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4303"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.94%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3746"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07339"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8769"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("D12").Value = "1.852.69"
$ws.Range("E12").Value = "  -0.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.738"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.433"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07130"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.41%  "

$ws.Range("E17").Value = "  -1.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008983"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("E19").Value = "  -1.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").Value = "27.716.61"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.204"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.82%  "

$ws.Range("E23").Value = "  -1.62%  "

$ws.Range("D24").Value = "2.076.60"
$ws.Range("E24").Value = "  -1.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.983"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.187"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.367"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08935"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.225"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7781"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.549"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.930"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.011"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.132"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01982"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("E39").Value = "  +0.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.276"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.898"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.15%  "

$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5128"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.819"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4774"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06469"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.011"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.687"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.839"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.67%  "
